$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value2 = 0
$ws.Range("J7").Value2 = 0
$ws.Range("L7").Value2 = 0
$ws.Range("N7").ClearContents()
$ws.Range("H13").Value2 = 57803.6
$ws.Range("J13").Value2 = 57803.6
$ws.Range("L13").Value2 = 57803.6
$ws.Range("N13").Value2 = -58141.6
$ws.Range("H14").Value2 = 0
$ws.Range("J14").Value2 = 0
$ws.Range("L14").Value2 = 0
$ws.Range("N14").ClearContents()
$ws.Range("H107").Value2 = 626
$ws.Range("I107").Value2 = 657.89746
$ws.Range("J107").Value2 = 552.82355
$ws.Range("K107").Value2 = 657.89746
$ws.Range("L107").Value2 = 552.82355
$ws.Range("M107").Value2 = 1262.10254
$ws.Range("N107").Value2 = -4392.82355
$ws.Range("H137").Value2 = 12502673
$ws.Range("I137").Value2 = 35716716
$ws.Range("J137").Value2 = 2803.8462
$ws.Range("K137").Value2 = 107150148
$ws.Range("L137").Value2 = 8411.5386
$ws.Range("M137").Value2 = -107147598
$ws.Range("N137").Value2 = -13511.5386

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 14498020
$ws.Range("I32").Value2 = 22224068
$ws.Range("J32").Value2 = 11679.375
$ws.Range("K32").Value2 = 22224068
$ws.Range("L32").Value2 = 11679.375
$ws.Range("M32").Value2 = -22223781
$ws.Range("N32").Value2 = -12253.375
$ws.Range("H33").Value2 = 0
$ws.Range("I33").Value2 = 0
$ws.Range("K33").Value2 = 0
$ws.Range("M33").ClearContents()
$ws.Range("H63").Value2 = 3111.111
$ws.Range("I63").Value2 = 2166.6667
$ws.Range("J63").Value2 = 5000
$ws.Range("K63").Value2 = 2166.6667
$ws.Range("L63").Value2 = 5000
$ws.Range("M63").Value2 = -1480.6667
$ws.Range("N63").Value2 = -6372
$ws.Range("H64").Value2 = 20000
$ws.Range("I64").Value2 = 0
$ws.Range("J64").Value2 = 20000
$ws.Range("K64").Value2 = 0
$ws.Range("L64").Value2 = 20000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value2 = -20496
$ws.Range("H66").Value2 = 3111.111
$ws.Range("I66").Value2 = 2166.6667
$ws.Range("J66").Value2 = 5000
$ws.Range("K66").Value2 = 10833.3335
$ws.Range("L66").Value2 = 25000
$ws.Range("M66").Value2 = -7401.333500000001
$ws.Range("N66").Value2 = -31864
$ws.Range("H67").Value2 = 20000
$ws.Range("I67").Value2 = 0
$ws.Range("J67").Value2 = 20000
$ws.Range("K67").Value2 = 0
$ws.Range("L67").Value2 = 20000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value2 = -21716
$ws.Range("H96").Value2 = 16511.6
$ws.Range("J96").Value2 = 16511.6
$ws.Range("L96").Value2 = 16511.6
$ws.Range("N96").Value2 = -22003.6

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1611.3572
$ws.Range("I31").Value2 = 1215.5385
$ws.Range("J31").Value2 = 1954.4
$ws.Range("K31").Value2 = 1215.5385
$ws.Range("L31").Value2 = 1954.4
$ws.Range("M31").Value2 = -920.5385000000001
$ws.Range("N31").Value2 = -2544.4
$ws.Range("H34").Value2 = 1611.3572
$ws.Range("I34").Value2 = 1215.5385
$ws.Range("J34").Value2 = 1954.4
$ws.Range("K34").Value2 = 1215.5385
$ws.Range("L34").Value2 = 1954.4
$ws.Range("M34").Value2 = -1013.5385
$ws.Range("N34").Value2 = -2358.4
$ws.Range("H56").Value2 = 8878.6
$ws.Range("I56").Value2 = 8131
$ws.Range("K56").Value2 = 8131
$ws.Range("M56").Value2 = -7286

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value2 = 1528.9
$ws.Range("I80").Value2 = 947.25
$ws.Range("J80").Value2 = 1916.6666
$ws.Range("K80").Value2 = 2841.75
$ws.Range("L80").Value2 = 5749.9998
$ws.Range("M80").Value2 = -1905.75
$ws.Range("N80").Value2 = -7621.9998
$ws.Range("H83").Value2 = 1528.9
$ws.Range("I83").Value2 = 947.25
$ws.Range("J83").Value2 = 1916.6666
$ws.Range("K83").Value2 = 8525.25
$ws.Range("L83").Value2 = 17249.9994
$ws.Range("M83").Value2 = -3845.25
$ws.Range("N83").Value2 = -26609.9994
$ws.Range("H101").Value2 = 4720
$ws.Range("J101").Value2 = 4950
$ws.Range("L101").Value2 = 14850
$ws.Range("N101").Value2 = -19718

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value2 = 8000
$ws.Range("J2").Value2 = 8000
$ws.Range("L2").Value2 = 8000
$ws.Range("N2").Value2 = -8224
$ws.Range("H22").Value2 = 501253.1
$ws.Range("I22").Value2 = 556786.75
$ws.Range("K22").Value2 = 556786.75
$ws.Range("M22").Value2 = -556491.75
$ws.Range("H26").Value2 = 0
$ws.Range("I26").Value2 = 0
$ws.Range("J26").Value2 = 0
$ws.Range("K26").Value2 = 0
$ws.Range("L26").Value2 = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value2 = 501253.1
$ws.Range("I27").Value2 = 556786.75
$ws.Range("K27").Value2 = 556786.75
$ws.Range("M27").Value2 = -556679.75
$ws.Range("H30").Value2 = 10450
$ws.Range("I30").Value2 = 10450
$ws.Range("J30").Value2 = 0
$ws.Range("K30").Value2 = 10450
$ws.Range("L30").Value2 = 0
$ws.Range("M30").Value2 = -10342
$ws.Range("N30").ClearContents()
$ws.Range("H61").Value2 = 7433.25
$ws.Range("I61").Value2 = 7819.9
$ws.Range("J61").Value2 = 5500
$ws.Range("K61").Value2 = 7819.9
$ws.Range("L61").Value2 = 5500
$ws.Range("M61").Value2 = -7617.9
$ws.Range("N61").Value2 = -5904
$ws.Range("H82").Value2 = 2869
$ws.Range("I82").Value2 = 2915.1428
$ws.Range("J82").Value2 = 2815.1667
$ws.Range("K82").Value2 = 2915.1428
$ws.Range("L82").Value2 = 2815.1667
$ws.Range("M82").Value2 = -2554.1428
$ws.Range("N82").Value2 = -3537.1667
$ws.Range("H85").Value2 = 2869
$ws.Range("I85").Value2 = 2915.1428
$ws.Range("J85").Value2 = 2815.1667
$ws.Range("K85").Value2 = 2915.1428
$ws.Range("L85").Value2 = 2815.1667
$ws.Range("M85").Value2 = -1667.1428
$ws.Range("N85").Value2 = -5311.1667
$ws.Range("H113").Value2 = 7433.25
$ws.Range("I113").Value2 = 7819.9
$ws.Range("J113").Value2 = 5500
$ws.Range("K113").Value2 = 7819.9
$ws.Range("L113").Value2 = 5500
$ws.Range("M113").Value2 = -5649.9
$ws.Range("N113").Value2 = -9840
